# Apply the "PO Forecast" update:
#  - rename the "Requested quantity" header on the existing sheets
#  - add a new "PO Forecast" sheet at the end with forecast data

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet so it lands at the end
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "PO Forecast"

# Header row
$ws3.Cells.Item(1, 1).Value = "ds"
$ws3.Cells.Item(1, 2).Value = "PO_Forecast"
$ws3.Cells.Item(1, 3).Value = "yhat_lower"
$ws3.Cells.Item(1, 4).Value = "yhat_upper"

$headerRng = $ws3.Range("A1:D1")
$headerRng.Font.Bold = $true
$headerRng.Borders.LineStyle = 1
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160

# Forecast rows: ds, PO_Forecast, yhat_lower, yhat_upper
$rows = @(
  @(45564.99999999999, 718, 249.9481739918187, 1161.111317244128),
  @(45599.99999999999, 560, 111.2698168678371, 1014.557140274866),
  @(45634.99999999999, 402, -54.04518084447861, 857.4982147038353),
  @(45641.99999999999, 370, -68.96576732883146, 857.4310544907519),
  @(45648.99999999999, 338, -126.5983080420662, 765.6081886421625),
  @(45655.99999999999, 307, -161.0090794016848, 753.9386037319946),
  @(45662.99999999999, 275, -180.1278112541868, 735.977849146814),
  @(45669.99999999999, 243, -225.2309538927573, 710.5668445471206),
  @(45676.99999999999, 212, -238.6375736994737, 690.2579157486206),
  @(45683.99999999999, 180, -277.1867163247519, 657.2641730656291),
  @(45690.99999999999, 149, -326.8695768554792, 611.4715420097945)
)

$r = 2
foreach ($row in $rows) {
  $ws3.Cells.Item($r, 1).Value = $row[0]
  $ws3.Cells.Item($r, 2).Value = $row[1]
  $ws3.Cells.Item($r, 3).Value = $row[2]
  $ws3.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

$dateRng = $ws3.Range("A2:A12")
$dateRng.NumberFormat = "YYYY-MM-DD HH:MM:SS"
